$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$timeTaken = @(
    "2021-10-05 10:50:28.709997",
    "2021-10-05 10:50:28.710008",
    "2021-10-05 10:50:28.710011",
    "2021-10-05 10:50:28.710014",
    "2021-10-05 10:50:28.710016",
    "2021-10-05 10:50:28.710019",
    "2021-10-05 10:50:28.710022",
    "2021-10-05 10:50:28.710024",
    "2021-10-05 10:50:28.710027",
    "2021-10-05 10:50:28.710029",
    "2021-10-05 10:50:28.710032",
    "2021-10-05 10:50:28.710034",
    "2021-10-05 10:50:28.710037",
    "2021-10-05 10:50:28.710039",
    "2021-10-05 10:50:28.710042",
    "2021-10-05 10:50:28.710044",
    "2021-10-05 10:50:28.710047",
    "2021-10-05 10:50:28.710050",
    "2021-10-05 10:50:28.710053",
    "2021-10-05 10:50:28.710055",
    "2021-10-05 10:50:28.710058",
    "2021-10-05 10:50:28.710060",
    "2021-10-05 10:50:28.710062",
    "2021-10-05 10:50:28.710065",
    "2021-10-05 10:50:28.710068",
    "2021-10-05 10:50:28.710071",
    "2021-10-05 10:50:28.710073",
    "2021-10-05 10:50:28.710076",
    "2021-10-05 10:50:28.710078",
    "2021-10-05 10:50:28.710080",
    "2021-10-05 10:50:28.710083",
    "2021-10-05 10:50:28.710085",
    "2021-10-05 10:50:28.710088",
    "2021-10-05 10:50:28.710091",
    "2021-10-05 10:50:28.710093",
    "2021-10-05 10:50:28.710096",
    "2021-10-05 10:50:28.710098",
    "2021-10-05 10:50:28.710101",
    "2021-10-05 10:50:28.710103",
    "2021-10-05 10:50:28.710106",
    "2021-10-05 10:50:28.710109",
    "2021-10-05 10:50:28.710112",
    "2021-10-05 10:50:28.710114",
    "2021-10-05 10:50:28.710116",
    "2021-10-05 10:50:28.710119",
    "2021-10-05 10:50:28.710121",
    "2021-10-05 10:50:28.710124",
    "2021-10-05 10:50:28.710127",
    "2021-10-05 10:50:28.710129",
    "2021-10-05 10:50:28.710132",
    "2021-10-05 10:50:28.710134",
    "2021-10-05 10:50:28.710137",
    "2021-10-05 10:50:28.710139",
    "2021-10-05 10:50:28.710142",
    "2021-10-05 10:50:28.710145",
    "2021-10-05 10:50:28.710147",
    "2021-10-05 10:50:28.710150",
    "2021-10-05 10:50:28.710152",
    "2021-10-05 10:50:28.710155",
    "2021-10-05 10:50:28.710157",
    "2021-10-05 10:50:28.710160",
    "2021-10-05 10:50:28.710162",
    "2021-10-05 10:50:28.710165",
    "2021-10-05 10:50:28.710167",
    "2021-10-05 10:50:28.710171",
    "2021-10-05 10:50:28.710174",
    "2021-10-05 10:50:28.710176",
    "2021-10-05 10:50:28.710179",
    "2021-10-05 10:50:28.710181",
    "2021-10-05 10:50:28.710184",
    "2021-10-05 10:50:28.710186",
    "2021-10-05 10:50:28.710189",
    "2021-10-05 10:50:28.710191",
    "2021-10-05 10:50:28.710193",
    "2021-10-05 10:50:28.710196",
    "2021-10-05 10:50:28.710199",
    "2021-10-05 10:50:28.710203",
    "2021-10-05 10:50:28.710206",
    "2021-10-05 10:50:28.710209",
    "2021-10-05 10:50:28.710211",
    "2021-10-05 10:50:28.710214",
    "2021-10-05 10:50:28.710216",
    "2021-10-05 10:50:28.710219",
    "2021-10-05 10:50:28.710221",
    "2021-10-05 10:50:28.710223",
    "2021-10-05 10:50:28.710226",
    "2021-10-05 10:50:28.710228",
    "2021-10-05 10:50:28.710231",
    "2021-10-05 10:50:28.710233",
    "2021-10-05 10:50:28.710236",
    "2021-10-05 10:50:28.710238",
    "2021-10-05 10:50:28.710241",
    "2021-10-05 10:50:28.710245",
    "2021-10-05 10:50:28.710248",
    "2021-10-05 10:50:28.710250",
    "2021-10-05 10:50:28.710253",
    "2021-10-05 10:50:28.710255",
    "2021-10-05 10:50:28.710258",
    "2021-10-05 10:50:28.710260",
    "2021-10-05 10:50:28.710263",
    "2021-10-05 10:50:28.710265",
    "2021-10-05 10:50:28.710268",
    "2021-10-05 10:50:28.710270",
    "2021-10-05 10:50:28.710273",
    "2021-10-05 10:50:28.710275",
    "2021-10-05 10:50:28.710278",
    "2021-10-05 10:50:28.710280",
    "2021-10-05 10:50:28.710283",
    "2021-10-05 10:50:28.710288",
    "2021-10-05 10:50:28.710291",
    "2021-10-05 10:50:28.710294",
    "2021-10-05 10:50:28.710296",
    "2021-10-05 10:50:28.710299",
    "2021-10-05 10:50:28.710301",
    "2021-10-05 10:50:28.710304",
    "2021-10-05 10:50:28.710306",
    "2021-10-05 10:50:28.710309",
    "2021-10-05 10:50:28.710312",
    "2021-10-05 10:50:28.710314",
    "2021-10-05 10:50:28.710317",
    "2021-10-05 10:50:28.710319",
    "2021-10-05 10:50:28.710322",
    "2021-10-05 10:50:28.710325",
    "2021-10-05 10:50:28.710327",
    "2021-10-05 10:50:28.710330",
    "2021-10-05 10:50:28.710332",
    "2021-10-05 10:50:28.710335",
    "2021-10-05 10:50:28.710337",
    "2021-10-05 10:50:28.710342",
    "2021-10-05 10:50:28.710345",
    "2021-10-05 10:50:28.710348",
    "2021-10-05 10:50:28.710351",
    "2021-10-05 10:50:28.710353",
    "2021-10-05 10:50:28.710356",
    "2021-10-05 10:50:28.710359",
    "2021-10-05 10:50:28.710361",
    "2021-10-05 10:50:28.710364",
    "2021-10-05 10:50:28.710366",
    "2021-10-05 10:50:28.710369",
    "2021-10-05 10:50:28.710372",
    "2021-10-05 10:50:28.710374",
    "2021-10-05 10:50:28.710377",
    "2021-10-05 10:50:28.710379",
    "2021-10-05 10:50:28.710382"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
